# Fruta / hortaliza, semanal
# Insert a new weekly record row right before the existing row 49
# (shifting all following rows down by one) and populate it with the
# new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 49; this pushes the former rows
# 49..85 down to 50..86 and extends the sheet dimension automatically.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly data point.
# Most columns mirror the record that used to sit at row 49 (same
# market/region/product grouping); only the date, quality grade,
# volume and price fields differ.
$ws.Range("A49").Value = 3
$ws.Range("B49").Value = "Femacal de La Calera"
$ws.Range("C49").Value = "Coquimbo"
$ws.Range("D49").Value = 44574
$ws.Range("E49").Value = 5
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100107
$ws.Range("H49").Value = "Otros"
$ws.Range("I49").Value = 100107011
$ws.Range("J49").Value = "Tuna"
$ws.Range("K49").Value = "Sin especificar"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 56
$ws.Range("N49").Value = 20000
$ws.Range("O49").Value = 20000
$ws.Range("P49").Value = 20000
$ws.Range("Q49").Value = "$/caja 16 kilos"
$ws.Range("R49").Value = "Provincia de Limarí"
$ws.Range("S49").Value = 1250
$ws.Range("T49").Value = 16

# Match the date cell's number format style used by the rest of column D.
$ws.Range("D49").NumberFormat = $ws.Range("D50").NumberFormat
